$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before current row 667, pushing the existing
# block (old rows 667-699) down to 669-701.
$ws.Rows.Item(667).EntireRow.Insert()
$ws.Rows.Item(667).EntireRow.Insert()

# New row 667 - "Primera" quality entry dated 44753
$ws.Cells.Item(667, 1).Value = 3
$ws.Cells.Item(667, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(667, 3).Value = "Coquimbo"
$ws.Cells.Item(667, 4).Value = 44753
$ws.Cells.Item(667, 5).Value = 5
$ws.Cells.Item(667, 6).Value = 100112008
$ws.Cells.Item(667, 7).Value = "Coliflor"
$ws.Cells.Item(667, 8).Value = "Sin especificar"
$ws.Cells.Item(667, 9).Value = "Primera"
$ws.Cells.Item(667, 10).Value = 1900
$ws.Cells.Item(667, 11).Value = 1000
$ws.Cells.Item(667, 12).Value = 1100
$ws.Cells.Item(667, 13).Value = 1050
$ws.Cells.Item(667, 14).Value = "`$/unidad"
$ws.Cells.Item(667, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(667, 16).Value = 1050
$ws.Cells.Item(667, 17).Value = 1
$ws.Cells.Item(667, 18).Value = "Hortaliza"

# New row 668 - "Segunda" quality entry dated 44753
$ws.Cells.Item(668, 1).Value = 3
$ws.Cells.Item(668, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(668, 3).Value = "Coquimbo"
$ws.Cells.Item(668, 4).Value = 44753
$ws.Cells.Item(668, 5).Value = 5
$ws.Cells.Item(668, 6).Value = 100112008
$ws.Cells.Item(668, 7).Value = "Coliflor"
$ws.Cells.Item(668, 8).Value = "Sin especificar"
$ws.Cells.Item(668, 9).Value = "Segunda"
$ws.Cells.Item(668, 10).Value = 900
$ws.Cells.Item(668, 11).Value = 800
$ws.Cells.Item(668, 12).Value = 800
$ws.Cells.Item(668, 13).Value = 800
$ws.Cells.Item(668, 14).Value = "`$/unidad"
$ws.Cells.Item(668, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(668, 16).Value = 800
$ws.Cells.Item(668, 17).Value = 1
$ws.Cells.Item(668, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D
$ws.Range("D667:D668").NumberFormat = $ws.Range("D669").NumberFormat
